$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 36 (shifts existing rows 36-39 down to 37-40)
$ws.Rows(36).Insert()

# Populate the new row 36 with the new weekly price record
$ws.Cells.Item(36, 1).Value2 = 1
$ws.Cells.Item(36, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(36, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(36, 4).Value2 = 44858
$ws.Cells.Item(36, 5).Value2 = 15
$ws.Cells.Item(36, 6).Value2 = 100112045
$ws.Cells.Item(36, 7).Value2 = "Zapallo"
$ws.Cells.Item(36, 8).Value2 = "Camote"
$ws.Cells.Item(36, 9).Value2 = "1a nueva(o)"
$ws.Cells.Item(36, 10).Value2 = 800
$ws.Cells.Item(36, 11).Value2 = 900
$ws.Cells.Item(36, 12).Value2 = 920
$ws.Cells.Item(36, 13).Value2 = 910
$ws.Cells.Item(36, 14).Value2 = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(36, 15).Value2 = "Perú"
$ws.Cells.Item(36, 16).Value2 = 910
$ws.Cells.Item(36, 17).Value2 = 1
$ws.Cells.Item(36, 18).Value2 = "Hortaliza"
